$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D
$ws.Range("D1").Value = "ITI"

# Updated data (Trial, Question, ConditionType, ITI) for rows 2-17
$data = @(
    @(1, 13, 1, 9),
    @(2, 21, 2, 6),
    @(3, 35, 2, 7),
    @(4, 31, 4, 7),
    @(5, 37, 4, 8),
    @(6, 1, 3, 6),
    @(7, 38, 4, 9),
    @(8, 2, 3, 7),
    @(9, 29, 2, 7),
    @(10, 23, 1, 8),
    @(11, 20, 1, 6),
    @(12, 26, 3, 7),
    @(13, 9, 3, 7),
    @(14, 33, 1, 6),
    @(15, 32, 4, 6),
    @(16, 4, 2, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Remove old rows 18, 19, 20 (trials 17, 18, 19) which no longer exist
$ws.Range("A18:D20").Clear()

# Update selection to match target
$ws.Range("D19").Select()
